$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the old "_GoBack" bookmark that sits after "When you
#    select a marker, it" (it will be re-created near the end of the
#    document, mirroring a real editing session where Word moves the
#    _GoBack bookmark to the location of the most recent edit).
# ------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

# ------------------------------------------------------------------
# 2) Fix the indentation of the "4.0  An accident marker ..." item so
#    it uses the same hanging-indent pattern as the other list items.
# ------------------------------------------------------------------
$rngFind = $d.Content
$rngFind.Find.Execute("An accident marker will not be visible") | Out-Null
$pAccident = $rngFind.Paragraphs.First
$pAccident.LeftIndent = 72
$pAccident.FirstLineIndent = -36

# ------------------------------------------------------------------
# 3) Add "...,  as accidents don't generally last longer" right before
#    the final period of that same sentence.
# ------------------------------------------------------------------
$rngFind = $d.Content
$rngFind.Find.Execute("has passed.") | Out-Null
$insertPoint = $d.Range($rngFind.End - 1, $rngFind.End - 1)
$newText1 = ", as accidents don" + [char]0x2019 + "t generally last longer"
$insertPoint.InsertBefore($newText1)

# keep the newly inserted text (and the trailing period that follows
# it) as their own runs instead of letting them silently merge back
# into the neighbouring identically-formatted runs
$startNew = $rngFind.End - 1
$endNew = $startNew + $newText1.Length
$rngNew = $d.Range($startNew, $endNew)
$rngNew.Font.Bold = $true
$rngNew.Font.Bold = $false
$rngPeriod = $d.Range($endNew, $endNew + 1)
$rngPeriod.Font.Bold = $true
$rngPeriod.Font.Bold = $false

# ------------------------------------------------------------------
# 4) Add ", or they can remove the marker altogether" right before the
#    final period of the "...turn the colour of the marker to green."
#    sentence.
# ------------------------------------------------------------------
$rngFind2 = $d.Content
$rngFind2.Find.Execute("turn the colour of the marker to green.") | Out-Null
$insertPoint2 = $d.Range($rngFind2.End - 1, $rngFind2.End - 1)
$newText2 = ", or they can remove the marker altogether"
$insertPoint2.InsertBefore($newText2)

$startNew2 = $rngFind2.End - 1
$endNew2 = $startNew2 + $newText2.Length
$rngNew2 = $d.Range($startNew2, $endNew2)
$rngNew2.Font.Bold = $true
$rngNew2.Font.Bold = $false
$rngPeriod2 = $d.Range($endNew2, $endNew2 + 1)
$rngPeriod2.Font.Bold = $true
$rngPeriod2.Font.Bold = $false

# ------------------------------------------------------------------
# 5) Collapse the four trailing empty, red (FF0000) "NormalWeb"
#    paragraphs at the end of the document down to a single one,
#    strip its left indent, and park the "_GoBack" bookmark inside it
#    (this is where Word leaves _GoBack after the edits above).
# ------------------------------------------------------------------
$redParas = @()
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Font.Color -eq 255) {
        $redParas += $i
    }
}

# The very last paragraph of the document body can't be removed by
# deleting its range (its paragraph mark is the document's final
# mark), so keep that one and delete the earlier, identical
# paragraphs instead, working from the highest index down so the
# indices we still have to visit stay valid.
for ($j = $redParas.Count - 2; $j -ge 0; $j--) {
    $d.Paragraphs($redParas[$j]).Range.Delete()
}

$lastIndex = $d.Paragraphs.Count
$pLast = $d.Paragraphs($lastIndex)
$pLast.LeftIndent = 0
$rLast = $pLast.Range
$rLast.Collapse(0)
$d.Bookmarks.Add("_GoBack", $rLast)
